$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Beverage sales section: add return functionality.
# Update Quantity (column E) for rows 2 and 3.
$ws.Range("E2").Value = 23
$ws.Range("E3").Value = 5

# Update Date (column G) and Time (column H) for row 3 (stored as plain text strings,
# not real date/time values). Temporarily force text format so Excel doesn't
# auto-convert these into date/time serial numbers, then restore the default style
# so the cells keep looking like ordinary (unformatted) text cells.
$ws.Range("G3:H3").NumberFormat = "@"
$ws.Range("G3").Value = "2024-09-15"
$ws.Range("H3").Value = "01:01:32"
$ws.Range("G3:H3").Style = "Normal"
